$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - HP
$ws.Range("D4").Value = 187
$ws.Range("E4").Value = 187
$ws.Range("H4").Value = 179
$ws.Range("I4").Value = 179
$ws.Range("L4").Value = 164
$ws.Range("M4").Value = 164

# Row 5 - MP
$ws.Range("D5").Value = 26
$ws.Range("E5").Value = 26
$ws.Range("H5").Value = 30
$ws.Range("I5").Value = 30
$ws.Range("L5").Value = 30
$ws.Range("M5").Value = 30

# Row 6 - Atk
$ws.Range("D6").Value = 24
$ws.Range("H6").Value = 10
$ws.Range("L6").Value = 8

# Row 7 - Def
$ws.Range("E7").Value = 0

# Row 8 - Level
$ws.Range("D8").Value = 10
$ws.Range("H8").Value = 10
$ws.Range("L8").Value = 10

# Row 9 - PWR
$ws.Range("D9").Value = 17
$ws.Range("H9").Value = 4
$ws.Range("L9").Value = 4

# Row 11 - HIT
$ws.Range("D11").Value = 10
$ws.Range("H11").Value = 11
$ws.Range("L11").Value = 10

# Row 12 - EV
$ws.Range("D12").Value = 11
$ws.Range("H12").Value = 9
$ws.Range("L12").Value = 8

# Row 13 - STM
$ws.Range("D13").Value = 22
$ws.Range("H13").Value = 13
$ws.Range("L13").Value = 11

# Row 14 - MAG
$ws.Range("D14").Value = 8
$ws.Range("H14").Value = 17
$ws.Range("L14").Value = 15

# Row 15 - MDEF
$ws.Range("D15").Value = 16
$ws.Range("E15").Value = 0
$ws.Range("H15").Value = 22
$ws.Range("L15").Value = 21
